# Implemented support for Docker, Zalenium and more optimization of codes
#
# RunManager sheet: remove the two obsolete test-case rows
# (checkWhetherGlobalViewIsSelectedByDefault, validateEnteringShipmentID),
# which shifts test1/test2/test3 up into rows 2-4. Also narrow the first two
# columns and move the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RunManager")

# Delete row 2 (checkWhetherGlobalViewIsSelectedByDefault) then the new row 2
# (validateEnteringShipmentID, shifted up after the first delete) so the
# remaining test1/test2/test3 rows move up to rows 2-4.
$ws.Rows.Item(2).Delete() | Out-Null
$ws.Rows.Item(2).Delete() | Out-Null

# Narrower columns for the trimmed-down sheet.
$ws.Columns.Item(1).ColumnWidth = 13.0
$ws.Columns.Item(2).ColumnWidth = 27.8

# Move the active selection.
$ws.Range("B4").Select() | Out-Null
